$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE) - update values
$ws.Range("B2").Value = 1.691
$ws.Range("C2").Value = 1.378
$ws.Range("D2").Value = 1.59
$ws.Range("E2").Value = 0.93
$ws.Range("F2").Value = 1.863

# Row 3 (MSE) - update values
$ws.Range("B3").Value = 4.178
$ws.Range("C3").Value = 3.963
$ws.Range("D3").Value = 4.426
$ws.Range("E3").Value = 1.442
$ws.Range("F3").Value = 5.451

# Row 4 previously held "mean Y-Test"; it becomes the "R2" row with new values
$ws.Range("A4").Value = "R2"
$ws.Range("B4").Value = 0.652
$ws.Range("C4").Value = 0.752
$ws.Range("D4").Value = 0.812
$ws.Range("E4").Value = 0.57
$ws.Range("F4").Value = 0.861

# Rows 5 (mean Y-predicted) and 6 (old R2) are removed entirely
$ws.Rows("5:6").Delete()
